# Auto-generated edits applying the diff for Jogos_da_Semana_FlashScore_2025-03-04.xlsx
# Updates odds values on Sheet1 across rows 5,6,15-20,22-25,27-29,32-34
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("AM5").Value = 7.5
$ws.Range("AP5").Value = 55

# Row 6
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 7

# Row 15
$ws.Range("G15").Value = 1.91

# Row 16
$ws.Range("I16").Value = 1.91
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10

# Row 17
$ws.Range("G17").Value = 2.88
$ws.Range("I17").Value = 2.45

# Row 18
$ws.Range("G18").Value = 1.75

# Row 19
$ws.Range("G19").Value = 2.5
$ws.Range("I19").Value = 2.8
$ws.Range("N19").Value = 9
$ws.Range("O19").Value = 1.4
$ws.Range("P19").Value = 2.75
$ws.Range("S19").Value = 2.25
$ws.Range("T19").Value = 1.62

# Row 20
$ws.Range("G20").Value = 1.7

# Row 22
$ws.Range("S22").Value = 2.2
$ws.Range("T22").Value = 1.65
$ws.Range("W22").Value = 4
$ws.Range("X22").Value = 1.22

# Row 23
$ws.Range("G23").Value = 2.38
$ws.Range("I23").Value = 3.1
$ws.Range("J23").Value = 3.1
$ws.Range("L23").Value = 3.6
$ws.Range("N23").Value = 9.5
$ws.Range("O23").Value = 1.3
$ws.Range("P23").Value = 3.4
$ws.Range("T23").Value = 1.8
$ws.Range("AD23").Value = 11
$ws.Range("AE23").Value = 9.5
$ws.Range("AG23").Value = 19
$ws.Range("AI23").Value = 9.5
$ws.Range("AM23").Value = 9.5
$ws.Range("AP23").Value = 34

# Row 24
$ws.Range("N24").Value = 9

# Row 25
$ws.Range("G25").Value = 1.8
$ws.Range("H25").Value = 3.8
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 2.4
$ws.Range("AD25").Value = 9.5
$ws.Range("AF25").Value = 15
$ws.Range("AG25").Value = 13
$ws.Range("AH25").Value = 23

# Row 27
$ws.Range("G27").Value = 1.33
$ws.Range("H27").Value = 5
$ws.Range("L27").Value = 7
$ws.Range("S27").Value = 1.53
$ws.Range("T27").Value = 2.4
$ws.Range("U27").Value = 1.9
$ws.Range("V27").Value = 1.9
$ws.Range("AA27").Value = 1.83
$ws.Range("AB27").Value = 1.83
$ws.Range("AD27").Value = 7
$ws.Range("AE27").Value = 9
$ws.Range("AF27").Value = 9
$ws.Range("AJ27").Value = 10
$ws.Range("AK27").Value = 21
$ws.Range("AM27").Value = 19

# Row 28
$ws.Range("G28").Value = 1.62
$ws.Range("H28").Value = 3.7
$ws.Range("I28").Value = 4.9
$ws.Range("J28").Value = 2.18
$ws.Range("L28").Value = 5.2
$ws.Range("O28").Value = 1.27
$ws.Range("P28").Value = 3.1
$ws.Range("S28").Value = 1.8
$ws.Range("T28").Value = 1.8
$ws.Range("W28").Value = 2.87
$ws.Range("X28").Value = 1.31
$ws.Range("Y28").Value = 1.39
$ws.Range("Z28").Value = 2.55
$ws.Range("AA28").Value = 1.82
$ws.Range("AB28").Value = 1.8
$ws.Range("AC28").Value = 6.7
$ws.Range("AD28").Value = 7.5
$ws.Range("AE28").Value = 8
$ws.Range("AF28").Value = 12
$ws.Range("AG28").Value = 13
$ws.Range("AI28").Value = 10
$ws.Range("AJ28").Value = 7.2
$ws.Range("AK28").Value = 16.5
$ws.Range("AL28").Value = 80
$ws.Range("AM28").Value = 13
$ws.Range("AN28").Value = 28
$ws.Range("AO28").Value = 16
$ws.Range("AP28").Value = 90
$ws.Range("AQ28").Value = 55
$ws.Range("AR28").Value = 55
$ws.Range("AS28").Value = 700

# Row 29
$ws.Range("H29").Value = 4.1
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 1.98
$ws.Range("K29").Value = 2.25
$ws.Range("O29").Value = 1.22
$ws.Range("P29").Value = 3.45
$ws.Range("S29").Value = 1.65
$ws.Range("T29").Value = 1.98
$ws.Range("W29").Value = 2.52
$ws.Range("X29").Value = 1.4
$ws.Range("AA29").Value = 1.78
$ws.Range("AB29").Value = 1.82
$ws.Range("AC29").Value = 7.3
$ws.Range("AD29").Value = 7.2
$ws.Range("AG29").Value = 11.5
$ws.Range("AH29").Value = 24
$ws.Range("AI29").Value = 12
$ws.Range("AJ29").Value = 8
$ws.Range("AM29").Value = 16.5
$ws.Range("AN29").Value = 37

# Row 32
$ws.Range("G32").Value = 3.75
$ws.Range("S32").Value = 1.83
$ws.Range("T32").Value = 2.03
$ws.Range("W32").Value = 3
$ws.Range("X32").Value = 1.36
$ws.Range("AC32").Value = 12
$ws.Range("AR32").Value = 26

# Row 33
$ws.Range("G33").Value = 3.5
$ws.Range("H33").Value = 3.3
$ws.Range("I33").Value = 2
$ws.Range("J33").Value = 4.33
$ws.Range("L33").Value = 2.75
$ws.Range("N33").Value = 10
$ws.Range("S33").Value = 2.08
$ws.Range("T33").Value = 1.73
$ws.Range("W33").Value = 3.75
$ws.Range("X33").Value = 1.25
$ws.Range("AA33").Value = 1.83
$ws.Range("AB33").Value = 1.83
$ws.Range("AC33").Value = 10
$ws.Range("AD33").Value = 17
$ws.Range("AE33").Value = 13
$ws.Range("AF33").Value = 41
$ws.Range("AG33").Value = 29
$ws.Range("AH33").Value = 41
$ws.Range("AI33").Value = 9
$ws.Range("AM33").Value = 7
$ws.Range("AN33").Value = 9
$ws.Range("AO33").Value = 9
$ws.Range("AP33").Value = 17
$ws.Range("AQ33").Value = 17
$ws.Range("AS33").Value = 301

# Row 34
$ws.Range("M34").Value = 1.07
$ws.Range("N34").Value = 9
$ws.Range("S34").Value = 2.3
$ws.Range("T34").Value = 1.6
